$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "New Holland Tractors"

# "4510" looks numeric, so force the cell to text format first so Excel
# stores it as a shared string (matching the original text-typed cell)
# instead of coercing it into a number.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "4510"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "['4510img0-4510-1632217675.png', '4510img1-upload-1632217675-0.png', '4510img2-4510-1632217675.png']"
